# #2 Implementada a Função calcular_cbr(resultados,cenarios)
# Adds a new "Parametros_Modular" worksheet (a modular copy of "Parametros")
# with crisis-factor columns, tweaks a couple of config/cost values, and
# nudges a few cell selections left over from the author's last session.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Configs: TaxaDeDesconto (F2) 0.1 -> 0, selection moves to B3
# ---------------------------------------------------------------------
$configs = $wb.Worksheets.Item("Configs")
$configs.Range("F2").Value = 0
$configs.Range("B3").Select()

# ---------------------------------------------------------------------
# 2. Create "Parametros_Modular" as a copy of "Parametros", placed right
#    before it (so sheet order becomes ... Cenarios, Parametros_Modular,
#    Parametros, Distribuições, ...)
# ---------------------------------------------------------------------
$parametros = $wb.Worksheets.Item("Parametros")
$parametros.Copy($parametros)
$modular = $wb.Worksheets.Item("Parametros (2)")
$modular.Name = "Parametros_Modular"

# New header cells I1:L1 - copy the bold header style from H1 then set text
$modular.Range("H1").Copy()
$modular.Range("I1").PasteSpecial(-4122)
$modular.Range("I1").Value = "Tem Crise?"

$modular.Range("H1").Copy()
$modular.Range("J1").PasteSpecial(-4122)
$modular.Range("J1").Value = "Fator Multiplicador"

$modular.Range("H1").Copy()
$modular.Range("K1").PasteSpecial(-4122)
$modular.Range("K1").Value = "Positivo?"

$modular.Range("H1").Copy()
$modular.Range("L1").PasteSpecial(-4122)
$modular.Range("L1").Value = "Impacto"

# Rows 2-5: crisis-factor inputs + parameter recalculated off Parametros
$modular.Range("C2").Formula = "=Parametros!C2*J2"
$modular.Range("I2").Formula = "=TRUE"
$modular.Range("J2").Formula = "=1+K2*L2"
$modular.Range("K2").Formula = "=TRUE"
$modular.Range("L2").Value = 0.1

$modular.Range("C3").Formula = "=Parametros!C3*J3"
$modular.Range("I3").Formula = "=TRUE"
$modular.Range("J3").Formula = "=1+K3*L3"
$modular.Range("K3").Formula = "=TRUE"
$modular.Range("L3").Value = 0.1

$modular.Range("C4").Formula = "=Parametros!C4*J4"
$modular.Range("I4").Formula = "=TRUE"
$modular.Range("J4").Formula = "=1+K4*L4"
$modular.Range("K4").Formula = "=TRUE"
$modular.Range("L4").Value = 0.1

$modular.Range("C5").Formula = "=Parametros!C5*J5"
$modular.Range("I5").Formula = "=TRUE"
$modular.Range("J5").Formula = "=1+K5*L5"
$modular.Range("K5").Formula = "=FALSE"
$modular.Range("L5").Value = 0.1

# Rows 6-17 (hidden scenario rows): just flag column I as TRUE
foreach ($r in 6..17) {
    $modular.Range("I$r").Formula = "=TRUE"
}

# Re-apply the AutoFilter on the copied range, now showing only
# "SemIniciativa" (matches the filtered state captured in the diff)
$modular.Range("A1:H33").AutoFilter(8, "SemIniciativa", 7)

# Register the (hidden) _FilterDatabase defined name for the new sheet
$modular.Names.Add("_xlnm._FilterDatabase", "=Parametros_Modular!`$A`$1:`$H`$33") | Out-Null
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Parametros_Modular!_FilterDatabase") {
        $n.Visible = $false
    }
}

# View tweaks matching the author's last session on this sheet
$modular.Application.ActiveWindow.Zoom = 115
$modular.Range("H39").Select()

# ---------------------------------------------------------------------
# 3. Parametros: values untouched, only the leftover selection changes
# ---------------------------------------------------------------------
$parametros.Range("F20").Select()

# ---------------------------------------------------------------------
# 4. Custos: rescale the benefit figures, selection moves to F10
# ---------------------------------------------------------------------
$custos = $wb.Worksheets.Item("Custos")
foreach ($r in 2..6) {
    $custos.Range("D$r").Value = 500
}
foreach ($r in 12..16) {
    $custos.Range("D$r").Value = 0
}
foreach ($r in 17..21) {
    $custos.Range("D$r").Value = 1500
}
$custos.Range("F10").Select()

# Make sure the copied/modular sheet ends up the active tab, mirroring
# the final state captured by the diff (activeTab points at it).
$modular.Activate()
